$d = $word.ActiveDocument

# --- 1. "Код:" section: strike whole paragraphs (incl. paragraph mark) ---
# "Add "converter" to views.py  (learn classes in Python)"
$d.Paragraphs(2).Range.Font.StrikeThrough = $true
# "Move css tokens from views.py to templates"
$d.Paragraphs(3).Range.Font.StrikeThrough = $true

# --- 2. "Модели:" section: strike whole paragraphs ---
# "Add "Cost" field to item"
$d.Paragraphs(6).Range.Font.StrikeThrough = $true
# "Improve Categories architecture"
$d.Paragraphs(7).Range.Font.StrikeThrough = $true

# --- 3. "Страниц «Контакты», «О компании», «Задать вопрос»" : strike only two words/phrases ---
$p13 = $d.Paragraphs(13)
$r = $p13.Range.Duplicate
$r.Find.Execute("Контакты")
$r.Font.StrikeThrough = $true

$r = $p13.Range.Duplicate
$r.Find.Execute("Задать вопрос")
$r.Font.StrikeThrough = $true

# --- 4. "Убрать «English Version»" : strike whole paragraph ---
$d.Paragraphs(14).Range.Font.StrikeThrough = $true

# --- 5. "Значки соц. сетей(Создание групп в соц. сетях)" : strike only "Значки соц. сетей" ---
$p15 = $d.Paragraphs(15)
$r = $p15.Range.Duplicate
$r.Find.Execute("Значки соц. сетей")
$r.Font.StrikeThrough = $true

# --- 6. "Масштабирование при сужении" : strike whole paragraph ---
$d.Paragraphs(17).Range.Font.StrikeThrough = $true

# --- 7. "Шрифт «кухня мебели»  в логотипе" : merge trailing runs + strike whole paragraph ---
$p20 = $d.Paragraphs(20)
$r = $p20.Range.Duplicate
$r.Find.Execute(" в логотипе")
Write-Output "found-v-logo:"
Write-Output $r.Find.Found
$p20.Range.Font.StrikeThrough = $true

# --- 8. "Белые поля" : strike whole paragraph ---
$d.Paragraphs(28).Range.Font.StrikeThrough = $true

# --- 9. "Когда одно фото, убрать кнопки" : strike whole paragraph ---
$d.Paragraphs(36).Range.Font.StrikeThrough = $true

# --- 10. "Оформить заказ" paragraph: merge runs, move bookmark, add page break paragraph ---
